# 5/19/20 update to Gantt chart
# - Adds a new "200515" sheet (copied/evolved from "200501") with the latest
#   task list, and makes small edits to the existing "200501" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the existing "200501" sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("200501")

# Category for the "Expanded dataset (12k points) exploration" row was
# renamed from "Expanded data" to "Chemical space".
$ws3.Range("B10").Value = "Chemical space"

# Column B got a bit wider to fit the new category names.
$ws3.Columns.Item(2).ColumnWidth = 13.14

# Move the current selection (cosmetic, but part of the saved view state).
$ws3.Range("F26").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) Create the new "200515" sheet as a copy of "200501" so it starts out
#    with identical formatting/styles, then bring its data up to date.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Copy([System.Reflection.Missing]::Value, $lastSheet) | Out-Null

$ws4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Name = "200515"

# Resize columns A and C to better fit the new (longer) task names.
$ws4.Columns.Item(1).ColumnWidth = 55.8
$ws4.Columns.Item(2).ColumnWidth = 13.14
$ws4.Columns.Item(3).ColumnWidth = 8.3

# Insert 4 additional rows to hold the new tasks added since 5/1/20.
# One before the old row 10 ("Expanded dataset..."), three more before the
# old row 12 ("Expected improvement...").
$ws4.Rows.Item(10).Insert()
$ws4.Rows.Item(12).Insert()
$ws4.Rows.Item(12).Insert()
$ws4.Rows.Item(12).Insert()

# --- Row 8: "Test models with mixed alloy data" dates slipped ---
$ws4.Range("C8").Value2 = 43973
$ws4.Range("D8").Value2 = 43987

# --- Row 10 (new): "Outlier detection and removal: PCA, KNN" ---
$ws4.Range("A10").Value = "Outlier detection and removal: PCA, KNN"
$ws4.Range("B10").Value = "Total data"
$ws4.Range("C10").Value2 = 43950
$ws4.Range("D10").Value2 = 43957

# --- Row 11: "Expanded dataset..." category renamed to "Chemical space" ---
$ws4.Range("B11").Value = "Chemical space"

# --- Row 12 (new): "Outlier detection and removal " ---
$ws4.Range("A12").Value = "Outlier detection and removal "
$ws4.Range("B12").Value = "Chemical space"
$ws4.Range("C12").Value2 = 43957
$ws4.Range("D12").Value2 = 43964

# --- Row 13 (new): "Uncertainty prediction: all models" ---
$ws4.Range("A13").Value = "Uncertainty prediction: all models"
$ws4.Range("B13").Value = "Chemical space"
$ws4.Range("C13").Value2 = 43957
$ws4.Range("D13").Value2 = 43971

# --- Row 14 (new): "Model tuning with new descriptors" ---
$ws4.Range("A14").Value = "Model tuning with new descriptors"
$ws4.Range("B14").Value = "New Total data"
$ws4.Range("C14").Value2 = 43966
$ws4.Range("D14").Value2 = 43978

# Rows 15-17 are the old rows 11-13 shifted down; their data is unchanged.

# Selection and active-sheet state for the new sheet.
$ws4.Range("D9").Select() | Out-Null
$ws4.Activate() | Out-Null
